# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename 'affiliations2' sheet to 'current_work' and populate it
# ------------------------------------------------------------------
$currentWork = $wb.Worksheets.Item("affiliations2")
$currentWork.Name = "current_work"

$currentWork.Range("A1").Value = "main"
$currentWork.Range("B1").Value = "sub"
$currentWork.Range("C1").Value = "working_with"
$currentWork.Range("D1").Value = "pic"
$currentWork.Range("E1").Value = "text"

$currentWork.Range("A2").Value = "ET"
$currentWork.Range("B2").Value = "Evapotranspiration Data Analysis"
$currentWork.Range("C2").Value = "Dami Eyelade, Jeremy Neill, Dino Korac"
$currentWork.Range("D2").Value = "resources/images/current_work/et.png"
$currentWork.Range("E2").Value = "Four projects look into modeling evapotranspiration across a continuous spatial surface within the contexts of the HydroData package. Specifically (1) how DAYMET and Weather Underground data can be used to calculate PET at a daily timesteps for North America; (2) how the same process can be applied globally; (4) how values – analogous to the FAO crop coefficients – can be developed for all vegetated Anderson level 2 land covers and (4) how to decouple the effect of temperature on ET rates using observed streamflow records."

$currentWork.Range("A3").Value = "Flood Mapping"
$currentWork.Range("C3").Value = "Jim Coll, Dinuke Munasinghe"
$currentWork.Range("D3").Value = "resources/images/current_work/floodmapping.gif"
$currentWork.Range("E3").Value = "Three projects look at the feasibility of real-time flood mapping using R and the National Water Model. The first (1) is to develop a system that allows GIS agencies throughout the country to install and run an hourly flood impacts models using the HAND methodology, the second (2) is to evaluate the accuracy of the HAND method against a repository of satellite derived flood extents and the third (3) is develop a filtering scheme to identify abnormally low and high flows in a NWM forecast allowing decision-makers and modelers to focus on areas that need attention."

$currentWork.Range("A4").Value = "Water Security"
$currentWork.Range("B4").Value = "Using Text-based Analysis to Simplify the Water Security Paradigm"
$currentWork.Range("C4").Value = "Keith Clarke"
$currentWork.Range("D4").Value = "resources/images/current_work/watersecurity.jpg"
$currentWork.Range("E4").Value = "Over the last 25 years the concept of water security has gained in popularity, but its meaning remains ambiguous despite several attempts to articulate a unified definition. This project seeks to better understand water security through a quantitative, text-based analysis of the literature to remove linguistic redundancies and overlap to identify core, timeless agreements."

$currentWork.Range("A5").Value = "Urban Growth Modeling"
$currentWork.Range("C5").Value = "Keith Clarke"
$currentWork.Range("D5").Value = "resources/images/current_work/urbangrowth.jpg"
$currentWork.Range("E5").Value = "This study looks to apply the SLEUTH land use/land cover model to the state of California to derive high resolution, probabilistic, long-range, Anderson level 2 land use products. Such products will be used for scenario-based hydroclimate studies with the NCAR Research Application Laboratory. "

$currentWork.Range("A6").Value = "Drought Models"
$currentWork.Range("B6").Value = "Agent-based Risk Models"
$currentWork.Range("C6").Value = "Marthe Wens, Jeroen Aerts, Ted Veldkamp"
$currentWork.Range("D6").Value = "resources/images/current_work/droughtmodels.jpg"
$currentWork.Range("E6").Value = "This work is part of a larger study being carried out at the Institute for Environmental Studies at VU Amsterdam looking to couple distributed hydrologic models and multi-actor, agent-based behavioral models within a drought risk context. One case study looks at applying this modeling strategy to a region in California’s Central Valley to better understand impacts on ground water withdrawals and aquifer levels."

# wrap text + column width for column E
$currentWork.Range("E1:E6").WrapText = $true

# row heights for wrapped rows
$currentWork.Rows.Item(2).RowHeight = 144
$currentWork.Rows.Item(3).RowHeight = 160
$currentWork.Rows.Item(4).RowHeight = 96
$currentWork.Rows.Item(5).RowHeight = 80
$currentWork.Rows.Item(6).RowHeight = 112

# column widths (values compensate for the engine's internal +5/6 padding)
$currentWork.Columns.Item(1).ColumnWidth = 19.998697916666668
$currentWork.Columns.Item(2).ColumnWidth = 54.166666666666664
$currentWork.Columns.Item(3).ColumnWidth = 35.166666666666664
$currentWork.Columns.Item(4).ColumnWidth = 41.166666666666664
$currentWork.Columns.Item(5).ColumnWidth = 54.166666666666664

$currentWork.Range("B2").Select() | Out-Null

# ------------------------------------------------------------------
# 2. Add new 'software' sheet at the end of the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$software = $wb.Worksheets.Add($null, $lastSheet)
$software.Name = "software"

$software.Range("A1").Value = "main"
$software.Range("B1").Value = "sub"
$software.Range("C1").Value = "working_with"
$software.Range("D1").Value = "pic"
$software.Range("E1").Value = "text"
$software.Range("F1").Value = "homepage"
$software.Range("G1").Value = "code"
$software.Range("H1").Value = "app"

$software.Range("A2").Value = "AOI"
$software.Range("B2").Value = "Area of Interest"
$software.Range("D2").Value = "resources/images/software/aoi.png"
$software.Range("E2").Value = "If you’ve ever found yourself needing to geocode or reverse geocode a location, formalize an area of interest, get bounding geometries, decribe a place by lat/long, or better understand spatial locations this package should be able to help.
An area of interest (AOI) is a geographic extent. It helps confine and formalize a unit of work to a geographic area, and prioritize and define research and sub setting efforts while improving reproducibility. They are built around concrete spatail attributes but often are discussed in a more colloquail way. The aim of the is package is to help make the colloquial understanding of space more concrete."
$software.Range("F2").Value = "https://mikejohnson51.github.io/AOI/"
$software.Range("G2").Value = "https://github.com/mikejohnson51/AOI"

$software.Range("A3").Value = "HydroData"
$software.Range("D3").Value = "resources/images/software/hydrodata.png"
$software.Range("E3").Value = "Almost all environmental research begins with data; the question of what data is available; and the challenge of gathering it. HydroData is an R package designed to help users find, get, visualize and use a range of Earth Systems data for a defined Area of Interest. The package provides functions to automatically download data from 19 sources; interactively visualize and share data within R and/or a browser; and to export data for external use in GIS."
$software.Range("F3").Value = "https://mikejohnson51.github.io/HydroData/"
$software.Range("G3").Value = "https://github.com/mikejohnson51/HydroData"

$software.Range("A4").Value = "FlowFinder"
$software.Range("B4").Value = "Hydrology as a Service"
$software.Range("C4").Value = "<a href=`"https://overdodactyl.github.io`">Pat Johnson</a>"
$software.Range("D4").Value = "resources/images/software/flowfinder.png"
$software.Range("E4").Value = "FlowFinder provides access to five-day out streamflow forecasts at a three hour time step for all of CONUS and visualizes data for a 225 square mile region centered on a users requested location (area of interest: AOI). Areas of exceedingly high flow - when compared to long term normals are marked and mapped throughout CONUS, and AOI metadata is provided. "
$software.Range("F4").Value = "https://mikejohnson51.github.io/FlowFinder/"
$software.Range("G4").Value = "https://github.com/mikejohnson51/FlowFinder"

$software.Range("A5").Value = "NWM"
$software.Range("B5").Value = "An R client for the National Water Model"
$software.Range("D5").Value = "resources/images/software/nwm.png"
$software.Range("E5").Value = "Each day the NOAA NWM produces ~400 GB of forecast data for the CONUS. This data is stored for a 40 day rolling window on the HydroShare Thredds server amounting to over 16,000 GB (2 TB) of data being stored and accessible at any one time. This package aims to provide access to this data in a clean, fast, and convenient way through the R environment."
$software.Range("F5").Value = "https://mikejohnson51.github.io/NWM/"
$software.Range("G5").Value = "ttps://github.com/mikejohnson51/NWM"

$software.Range("A6").Value = "nwmViewer"
$software.Range("B6").Value = "Vizualize National Water Model Output"
$software.Range("C6").Value = "<a href=`"https://overdodactyl.github.io`">Pat Johnson</a>"
$software.Range("D6").Value = "resources/images/software/nwmViewer.png"
$software.Range("E6").Value = "nwmViewer is an R Shiny application built to quickly find and vizualize data via the nwm package. "
$software.Range("G6").Value = "https://github.com/overdodactyl/nwmViewer"
$software.Range("H6").Value = "https://nwmviewer.shinyapps.io/view/"

# wrap text column E
$software.Range("E1:E6").WrapText = $true

# row heights for wrapped rows
$software.Rows.Item(2).RowHeight = 176
$software.Rows.Item(3).RowHeight = 128
$software.Rows.Item(4).RowHeight = 96
$software.Rows.Item(5).RowHeight = 96
$software.Rows.Item(6).RowHeight = 32

# column widths
$software.Columns.Item(2).ColumnWidth = 33.498697916666664
$software.Columns.Item(4).ColumnWidth = 36.666666666666664
$software.Columns.Item(5).ColumnWidth = 53.166666666666664
$software.Columns.Item(6).ColumnWidth = 37.166666666666664
$software.Columns.Item(7).ColumnWidth = 38.830729166666664
$software.Columns.Item(8).ColumnWidth = 31.998697916666668

$software.Range("I1").Select() | Out-Null
$software.Activate() | Out-Null

# ------------------------------------------------------------------
# 3. Update the 'affiliations' sheet selection (no longer the active tab)
# ------------------------------------------------------------------
$affiliations = $wb.Worksheets.Item("affiliations")
$affiliations.Range("B3").Select() | Out-Null

# ------------------------------------------------------------------
# 4. Workbook window view settings
# ------------------------------------------------------------------
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 620
    $win.Top = 0
    $win.Width = 25540
    $win.Height = 15560
} catch {}

try {
    $win = $wb.Windows.Item(1)
    $win.FirstVisibleSheet = $wb.Worksheets.Item("awards").Index
} catch {}

# software tab is the final activated/selected sheet -> becomes activeTab
$software.Activate() | Out-Null
